$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old content in A1 (no longer used)
$ws.Range("A1").ClearContents()

# Row 1 - first user record
$ws.Range("B1").Value = "Sangeethapriya"
$ws.Range("C1").Value = "P R"
$ws.Range("D1").Value = "stest@gmail.com"
$ws.Range("E1").Value = "Stest@123"
$ws.Range("F1").Value = "Stest@123"

# Row 2 - second user record
$ws.Range("B2").Value = "Priya"
$ws.Range("C2").Value = "S"
$ws.Range("D2").Value = "Priyas@gmail.com"
$ws.Range("E2").Value = "priyanew@123"
$ws.Range("F2").Value = "priyanew@123"

# Hyperlinks
$ws.Hyperlinks.Add($ws.Range("D1"), "stest@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E1"), "Stest@123")
$ws.Hyperlinks.Add($ws.Range("F1"), "Stest@123")
$ws.Hyperlinks.Add($ws.Range("D2"), "Priyas@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E2"), "priyanew@123")
$ws.Hyperlinks.Add($ws.Range("F2"), "priyanew@123")

$ws.Range("F2").Select()
